$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.905.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.793.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.80"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5403"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3795"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07428"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.80"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.087"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.200"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.409"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.29"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.797.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.63"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001059"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06495"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.924"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.935.73"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.16"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.33"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.000.73"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.336"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.99"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.106"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.660"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.523"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06937"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2205"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02281"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.047"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.478"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.34"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.417"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.163"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.32"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.680"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5698"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.42"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.175"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.909"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06794"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.64"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.42%  "
